$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now holds the Newmont Corporation (NEM) data
$ws.Range("B2").Value = "Newmont Corporation"
$ws.Range("C2").Value = "NEM"
$ws.Range("D2").Value = 90.73
$ws.Range("F2").Value = 10.97
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 69.90000000000001
$ws.Range("N2").Value = 85.77505782882612

# Row 3 now holds the StreetTRACKS Gold Shares (GLD) data
$ws.Range("B3").Value = "StreetTRACKS Gold Shares"
$ws.Range("C3").Value = "GLD"
$ws.Range("D3").Value = 387.88
$ws.Range("F3").Value = 3.48
$ws.Range("H3").Value = 56
$ws.Range("J3").Value = 73
$ws.Range("N3").Value = 85.77505782882612

# Row 4 (Gold Dec 25 futures) keeps its own data, only MACRO_SCORE refreshed
$ws.Range("N4").Value = 85.77505782882612
